# Duplicate the "Tbl_LOGICollectionsManagementR1" sheet to create a new
# "Tbl_LOGICollectionsManagementR2" sheet (same header row / layout), then
# repoint the workbook-level defined name
# "Tbl_LOGICollectionsManagementReport" at the new sheet.

$wb = $excel.ActiveWorkbook

# Source sheet to clone.
$sourceSheet = $wb.Worksheets.Item('Tbl_LOGICollectionsManagementR1')

# Copy it so it lands right after itself, then rename the copy.
$sourceSheet.Copy([System.Reflection.Missing]::Value, $sourceSheet)
$newSheet = $wb.Worksheets.Item($sourceSheet.Index + 1)
$newSheet.Name = 'Tbl_LOGICollectionsManagementR2'

# Point the defined name at the newly added sheet.
$wb.Names.Item('Tbl_LOGICollectionsManagementReport').RefersTo = '=Tbl_LOGICollectionsManagementR2!$A$1:$W$1'
